# Updated cryptos list values (Price / Volume(1h) columns) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.052.00'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '2.237.92'
$ws.Range("E3").Value = '  -3.06%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.46'
$ws.Range("E5").Value = '  -2.61%  '
$ws.Range("E6").Value = '  -1.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.77'
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.626'
$ws.Range("E9").Value = '  -4.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.31'
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0952'
$ws.Range("E11").Value = '  -4.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.11'
$ws.Range("E12").Value = '  -3.93%  '
$ws.Range("E13").Value = '  -2.03%  '
$ws.Range("D14").Value = '2.573.23'
$ws.Range("E14").Value = '  -3.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.85'
$ws.Range("E15").Value = '  -4.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.860'
$ws.Range("E16").Value = '  -2.53%  '
$ws.Range("D17").Value = '2.255.16'
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").Value = '41.901.13'
$ws.Range("E18").Value = '  -2.77%  '
$ws.Range("D19").Value = '0.0₃0979'
$ws.Range("E19").Value = '  -2.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.13'
$ws.Range("E20").Value = '  -3.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.62'
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.27'
$ws.Range("E22").Value = '  +2.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.37'
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.47'
$ws.Range("E24").Value = '  -1.69%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  -6.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.30'
$ws.Range("E27").Value = '  -5.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.19'
$ws.Range("E28").Value = '  +12.08%  '
$ws.Range("E29").Value = '  -1.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.05'
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.54'
$ws.Range("E31").Value = '  -3.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.42'
$ws.Range("E32").Value = '  +3.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0840'
$ws.Range("E33").Value = '  +2.60%  '
$ws.Range("E34").Value = '  -5.47%  '
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.52'
$ws.Range("E36").Value = '  -2.75%  '
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.36'
$ws.Range("E39").Value = '  -7.82%  '
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("E41").Value = '  -7.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '111.92'
$ws.Range("E42").Value = '  +13.02%  '
$ws.Range("E43").Value = '  -5.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.47'
$ws.Range("E44").Value = '  -3.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.71'
$ws.Range("E45").Value = '  -5.09%  '
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("E48").Value = '  -4.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.34'
$ws.Range("E49").Value = '  -11.61%  '
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.17'
$ws.Range("E51").Value = '  -2.89%  '
